$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.904.78'
$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").Value = '1.701.13'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.66%  '

$ws.Range("D5").Value = '''315.29'
$ws.Range("E5").Value = '  -0.06%  '

$ws.Range("E6").Value = '  -0.47%  '

$ws.Range("D7").Value = '''0.4015'
$ws.Range("E7").Value = '  +2.26%  '

$ws.Range("D8").Value = '''0.4067'
$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = '''1.002'
$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").Value = '''53.79'
$ws.Range("E10").Value = '  +1.17%  '

$ws.Range("E11").Value = '  -2.87%  '

$ws.Range("D12").Value = '''0.08821'
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("D13").Value = '''25.84'
$ws.Range("E13").Value = '  +6.18%  '

$ws.Range("D14").Value = '''7.485'
$ws.Range("E14").Value = '  -1.90%  '

$ws.Range("D15").Value = '''8.054'
$ws.Range("E15").Value = '  +0.76%  '

$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").Value = '1.793.47'
$ws.Range("E17").Value = '  +5.77%  '

$ws.Range("D18").Value = '''96.68'
$ws.Range("E18").Value = '  -2.03%  '

$ws.Range("D19").Value = '''0.07200'
$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("D20").Value = '''20.98'
$ws.Range("E20").Value = '  +5.94%  '

$ws.Range("D21").Value = '''7.248'
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").Value = '''1.003'

$ws.Range("D23").Value = '''14.55'
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("D24").Value = '24.905.82'
$ws.Range("E24").Value = '  +1.13%  '

$ws.Range("D25").Value = '''2.337'
$ws.Range("E25").Value = '  -1.06%  '

$ws.Range("D26").Value = '''2.893'
$ws.Range("E26").Value = '  -5.14%  '

$ws.Range("D27").Value = '''6.692'
$ws.Range("E27").Value = '  +28.08%  '

$ws.Range("E28").Value = '  +1.71%  '

$ws.Range("D29").Value = '''163.43'
$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("D30").Value = '''143.55'
$ws.Range("E30").Value = '  +4.08%  '

$ws.Range("D31").Value = '''8.201'
$ws.Range("E31").Value = '  -3.12%  '

$ws.Range("D32").Value = '1.948.66'
$ws.Range("E32").Value = '  +3.42%  '

$ws.Range("D33").Value = '''2.275'
$ws.Range("E33").Value = '  +14.43%  '

$ws.Range("D34").Value = '''0.08738'

$ws.Range("D35").Value = '''7.405'
$ws.Range("E35").Value = '  -0.54%  '

$ws.Range("D36").Value = '''0.03173'
$ws.Range("E36").Value = '  +8.56%  '

$ws.Range("D37").Value = '''1.039'
$ws.Range("E37").Value = '  -1.18%  '

$ws.Range("D38").Value = '''0.2876'
$ws.Range("E38").Value = '  +5.28%  '

$ws.Range("D39").Value = '''0.8548'
$ws.Range("E39").Value = '  +8.65%  '

$ws.Range("D40").Value = '''10.89'
$ws.Range("E40").Value = '  +0.77%  '

$ws.Range("D41").Value = '''0.09440'
$ws.Range("E41").Value = '  +3.20%  '

$ws.Range("D42").Value = '''14.08'
$ws.Range("E42").Value = '  -1.49%  '

$ws.Range("D43").Value = '''1.472'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").Value = '''17.58'
$ws.Range("E44").Value = '  +6.02%  '

$ws.Range("D45").Value = '''2.705'
$ws.Range("E45").Value = '  +5.47%  '

$ws.Range("D46").Value = '''0.7478'
$ws.Range("E46").Value = '  +3.73%  '

$ws.Range("D47").Value = '''4.230'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Value = '''1.410'
$ws.Range("E48").Value = '  +5.89%  '

$ws.Range("D49").Value = '''1.002'
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("D50").Value = '''141.33'
$ws.Range("E50").Value = '  +1.44%  '

$ws.Range("D51").Value = '''0.08350'
$ws.Range("E51").Value = '  +4.55%  '
